$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.855.78"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "1.755.37"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'326.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'0.4590"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'0.3493"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.71%  "
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("D10").Value = "'0.07347"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("D11").Value = "'1.078"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").Value = "'20.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.55%  "
$ws.Range("D14").Value = "'5.969"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("D15").Value = "'7.135"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "1.753.91"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "'91.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.77%  "
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Value = "'0.06407"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "'16.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.43%  "
$ws.Range("D22").Value = "'5.748"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").Value = "27.875.73"
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("D24").Value = "'11.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").Value = "'2.166"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.74%  "
$ws.Range("D26").Value = "'162.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.12%  "
$ws.Range("D27").Value = "'19.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").Value = "1.956.98"
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").Value = "'2.129"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "'122.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.05%  "
$ws.Range("D31").Value = "'1.062"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("D32").Value = "'0.09245"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("D33").Value = "'3.665"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").Value = "'5.520"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").Value = "'11.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D36").Value = "'0.02262"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.70%  "
$ws.Range("D37").Value = "'0.06068"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("D38").Value = "'0.2058"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("D39").Value = "'4.887"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("D40").Value = "'0.6145"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("D41").Value = "'1.177"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").Value = "'7.755"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").Value = "'1.350"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("D44").Value = "'13.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D45").Value = "'3.730"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("D46").Value = "'0.5757"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("D47").Value = "'122.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("D48").Value = "'1.920"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").Value = "'0.06798"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'1.117"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").Value = "'71.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.36%  "
